$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '63.005.35'
$ws.Range('D2').Style = "Normal"

$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -1.61%  '
$ws.Range('E2').Style = "Normal"

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.189.96'
$ws.Range('D3').Style = "Normal"

$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +1.58%  '
$ws.Range('E3').Style = "Normal"

$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('E4').Style = "Normal"

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '590.46'
$ws.Range('D5').Style = "Normal"

$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -2.01%  '
$ws.Range('E5').Style = "Normal"

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '139.68'
$ws.Range('D6').Style = "Normal"

$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -2.41%  '
$ws.Range('E6').Style = "Normal"

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.185.34'
$ws.Range('D8').Style = "Normal"

$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +1.57%  '
$ws.Range('E8').Style = "Normal"

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.519'
$ws.Range('D9').Style = "Normal"

$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -0.90%  '
$ws.Range('E9').Style = "Normal"

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.147'
$ws.Range('D10').Style = "Normal"

$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -1.46%  '
$ws.Range('E10').Style = "Normal"

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.40'
$ws.Range('D11').Style = "Normal"

$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +0.24%  '
$ws.Range('E11').Style = "Normal"

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.462'
$ws.Range('D12').Style = "Normal"

$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -1.11%  '
$ws.Range('E12').Style = "Normal"

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000247'
$ws.Range('D13').Style = "Normal"

$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -2.97%  '
$ws.Range('E13').Style = "Normal"

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '34.22'
$ws.Range('D14').Style = "Normal"

$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -2.41%  '
$ws.Range('E14').Style = "Normal"

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.705.06'
$ws.Range('D15').Style = "Normal"

$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +1.53%  '
$ws.Range('E15').Style = "Normal"

$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +0.70%  '
$ws.Range('E16').Style = "Normal"

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.174.38'
$ws.Range('D17').Style = "Normal"

$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +1.35%  '
$ws.Range('E17').Style = "Normal"

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '63.008.16'
$ws.Range('D18').Style = "Normal"

$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -1.54%  '
$ws.Range('E18').Style = "Normal"

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.72'
$ws.Range('D19').Style = "Normal"

$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -1.62%  '
$ws.Range('E19').Style = "Normal"

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '474.95'
$ws.Range('D20').Style = "Normal"

$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -2.52%  '
$ws.Range('E20').Style = "Normal"

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.03'
$ws.Range('D21').Style = "Normal"

$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -4.63%  '
$ws.Range('E21').Style = "Normal"

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.708'
$ws.Range('D22').Style = "Normal"

$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('E22').Style = "Normal"

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.80'
$ws.Range('D23').Style = "Normal"

$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +1.99%  '
$ws.Range('E23').Style = "Normal"

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '83.85'
$ws.Range('D24').Style = "Normal"

$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -3.60%  '
$ws.Range('E24').Style = "Normal"

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '13.09'
$ws.Range('D25').Style = "Normal"

$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -2.63%  '
$ws.Range('E25').Style = "Normal"

$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('E26').Style = "Normal"

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.72'
$ws.Range('D27').Style = "Normal"

$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -1.44%  '
$ws.Range('E27').Style = "Normal"

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.12'
$ws.Range('D28').Style = "Normal"

$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +1.87%  '
$ws.Range('E28').Style = "Normal"

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.00'
$ws.Range('D29').Style = "Normal"

$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -3.21%  '
$ws.Range('E29').Style = "Normal"

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.10'
$ws.Range('D30').Style = "Normal"

$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +1.65%  '
$ws.Range('E30').Style = "Normal"

$ws.Range('B31').NumberFormat = "@"
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('B31').Style = "Normal"

$ws.Range('C31').NumberFormat = "@"
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('C31').Style = "Normal"

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '27.07'
$ws.Range('D31').Style = "Normal"

$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -0.19%  '
$ws.Range('E31').Style = "Normal"

$ws.Range('B32').NumberFormat = "@"
$ws.Range('B32').Value = 'FirstDigitalUSD'
$ws.Range('B32').Style = "Normal"

$ws.Range('C32').NumberFormat = "@"
$ws.Range('C32').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('C32').Style = "Normal"

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.00'
$ws.Range('D32').Style = "Normal"

$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +0.05%  '
$ws.Range('E32').Style = "Normal"

$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -2.84%  '
$ws.Range('E33').Style = "Normal"

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.54'
$ws.Range('D34').Style = "Normal"

$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -4.00%  '
$ws.Range('E34').Style = "Normal"

$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -2.45%  '
$ws.Range('E35').Style = "Normal"

$ws.Range('B36').NumberFormat = "@"
$ws.Range('B36').Value = 'Filecoin'
$ws.Range('B36').Style = "Normal"

$ws.Range('C36').NumberFormat = "@"
$ws.Range('C36').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('C36').Style = "Normal"

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.82'
$ws.Range('D36').Style = "Normal"

$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -3.22%  '
$ws.Range('E36').Style = "Normal"

$ws.Range('B37').NumberFormat = "@"
$ws.Range('B37').Value = 'OKB'
$ws.Range('B37').Style = "Normal"

$ws.Range('C37').NumberFormat = "@"
$ws.Range('C37').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('C37').Style = "Normal"

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '52.61'
$ws.Range('D37').Style = "Normal"

$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +0.14%  '
$ws.Range('E37').Style = "Normal"

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0₃0714'
$ws.Range('D38').Style = "Normal"

$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -4.06%  '
$ws.Range('E38').Style = "Normal"

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0389'
$ws.Range('D39').Style = "Normal"

$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -1.58%  '
$ws.Range('E39').Style = "Normal"

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '422.23'
$ws.Range('D40').Style = "Normal"

$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -3.65%  '
$ws.Range('E40').Style = "Normal"

$ws.Range('B41').NumberFormat = "@"
$ws.Range('B41').Value = 'Cosmos'
$ws.Range('B41').Style = "Normal"

$ws.Range('C41').NumberFormat = "@"
$ws.Range('C41').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('C41').Style = "Normal"

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '8.38'
$ws.Range('D41').Style = "Normal"

$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +1.18%  '
$ws.Range('E41').Style = "Normal"

$ws.Range('B42').NumberFormat = "@"
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('B42').Style = "Normal"

$ws.Range('C42').NumberFormat = "@"
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('C42').Style = "Normal"

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.76'
$ws.Range('D42').Style = "Normal"

$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -7.21%  '
$ws.Range('E42').Style = "Normal"

$ws.Range('B43').NumberFormat = "@"
$ws.Range('B43').Value = 'Maker'
$ws.Range('B43').Style = "Normal"

$ws.Range('C43').NumberFormat = "@"
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('C43').Style = "Normal"

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.964.27'
$ws.Range('D43').Style = "Normal"

$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +2.97%  '
$ws.Range('E43').Style = "Normal"

$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -7.70%  '
$ws.Range('E44').Style = "Normal"

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.264'
$ws.Range('D45').Style = "Normal"

$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +1.54%  '
$ws.Range('E45').Style = "Normal"

$ws.Range('B46').NumberFormat = "@"
$ws.Range('B46').Value = 'Fetch.AI'
$ws.Range('B46').Style = "Normal"

$ws.Range('C46').NumberFormat = "@"
$ws.Range('C46').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('C46').Style = "Normal"

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.14'
$ws.Range('D46').Style = "Normal"

$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -2.93%  '
$ws.Range('E46').Style = "Normal"

$ws.Range('B47').NumberFormat = "@"
$ws.Range('B47').Value = 'USDe'
$ws.Range('B47').Style = "Normal"

$ws.Range('C47').NumberFormat = "@"
$ws.Range('C47').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('C47').Style = "Normal"

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.999'
$ws.Range('D47').Style = "Normal"

$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +0.04%  '
$ws.Range('E47').Style = "Normal"

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '25.73'
$ws.Range('D48').Style = "Normal"

$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -1.00%  '
$ws.Range('E48').Style = "Normal"

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.114'
$ws.Range('D49').Style = "Normal"

$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -0.17%  '
$ws.Range('E49').Style = "Normal"

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.26'
$ws.Range('D50').Style = "Normal"

$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -5.11%  '
$ws.Range('E50').Style = "Normal"

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '119.73'
$ws.Range('D51').Style = "Normal"

$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -1.07%  '
$ws.Range('E51').Style = "Normal"
